$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert 10 new rows starting at row 17 ---
# This pushes old row17 (blank) -> row27, ukupno(18)->28, uplaceno(19)->29, and the
# trailing blank rows (20-30) -> (30-40). New rows 17-26 are created blank.
$ws.Rows("17:26").Insert()

# --- 2. Fill the new data rows (17-25) ---
$newRows = @(
    @(17, 21.07, "50 min",  1000),
    @(18, 29.07, "245 min", 4900),
    @(19, 30.07, "230 min", 4600),
    @(20, 2.08,  "280 min", 5600),
    @(21, 3.08,  "305 min", 6010),
    @(22, 4.08,  "150 min", 3000),
    @(23, 6.08,  "345 min", 6900),
    @(24, 9.08,  "450 min", 9000),
    @(25, 10.08, "45 min",  900)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
}

# --- 3. Update the totals block (now at rows 28/29) ---
$ws.Range("D28").Formula = "=SUM(D4:D25)"
$ws.Range("D29").Value = 20000

# --- 4. Number formatting on the amount column (Accounting, no symbol) ---
$acctFmt = "_(* #,##0_);_(* \(#,##0\);_(* ""-""_);_(@_)"

$ws.Range("D4:D11").HorizontalAlignment = 1
$ws.Range("D4:D11").NumberFormat = $acctFmt

# D16 carried an explicit bottom border from before the insert; clear it first so
# D12:D16 end up uniformly borderless after the accounting format is applied.
$ws.Range("D16").Borders.Item(9).LineStyle = -4142
$ws.Range("D12:D16").HorizontalAlignment = 1
$ws.Range("D12:D16").NumberFormat = $acctFmt

$ws.Range("D17:D24").HorizontalAlignment = 1
$ws.Range("D17:D24").NumberFormat = $acctFmt

$ws.Range("D25").HorizontalAlignment = -4152
$ws.Range("D25").NumberFormat = $acctFmt

$ws.Range("D26").HorizontalAlignment = -4108
$ws.Range("D26").NumberFormat = $acctFmt

$ws.Range("D27").HorizontalAlignment = -4108
$ws.Range("D27").NumberFormat = $acctFmt

$ws.Range("D28:D29").HorizontalAlignment = -4152
$ws.Range("D28:D29").NumberFormat = $acctFmt

# --- 5. Right-align the date/category columns for all data rows (plus totals labels) ---
$ws.Range("B4:C25").HorizontalAlignment = -4152
$ws.Range("C28:C29").HorizontalAlignment = -4152

# --- 6. Header row: bold font + bottom border (already present) but right aligned ---
$ws.Range("B3:D3").HorizontalAlignment = -4152

# --- 7. Column D width (matches bestFit width picked up from the Accounting format) ---
$ws.Columns("D").ColumnWidth = 11.5703125

# --- 8. Selection / view state ---
$ws.Range("J16").Select()
